$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "62.909.00"
Set-TextValue "E2" "  -1.34%  "
Set-TextValue "D3" "2.542.61"
Set-TextValue "E3" "  -0.07%  "
Set-TextValue "E4" "  +0.03%  "
Set-TextValue "D5" "574.38"
Set-TextValue "E5" "  +0.05%  "
Set-TextValue "D6" "145.83"
Set-TextValue "E6" "  -1.40%  "
Set-TextValue "E7" "  +0.02%  "
Set-TextValue "E8" "  -1.20%  "
Set-TextValue "E9" "  -1.65%  "
Set-TextValue "D10" "5.49"
Set-TextValue "E10" "  -4.75%  "
Set-TextValue "E12" "  -1.08%  "
Set-TextValue "D13" "27.25"
Set-TextValue "E13" "  -3.16%  "
Set-TextValue "D14" "2.996.90"
Set-TextValue "E14" "  +0.00%  "
Set-TextValue "D15" "62.831.31"
Set-TextValue "E15" "  -1.14%  "
Set-TextValue "E16" "  -1.11%  "
Set-TextValue "D17" "2.542.51"
Set-TextValue "E17" "  -0.14%  "
Set-TextValue "E18" "  -2.08%  "
Set-TextValue "D19" "335.91"
Set-TextValue "E19" "  -1.63%  "
Set-TextValue "D20" "4.31"
Set-TextValue "E20" "  -0.92%  "
Set-TextValue "E21" "  -2.18%  "
Set-TextValue "E22" "  +0.12%  "
Set-TextValue "D23" "65.21"
Set-TextValue "E23" "  -1.51%  "
Set-TextValue "E24" "  -0.34%  "
Set-TextValue "E25" "  +0.87%  "
Set-TextValue "E26" "  -0.03%  "
Set-TextValue "D27" "8.32"
Set-TextValue "E27" "  -0.09%  "
Set-TextValue "E28" "  +1.60%  "
Set-TextValue "D29" "7.23"
Set-TextValue "E29" "  +4.22%  "
Set-TextValue "E30" "  -2.80%  "
Set-TextValue "E31" "  -1.22%  "
Set-TextValue "D32" "177.63"
Set-TextValue "E32" "  +0.48%  "
Set-TextValue "E33" "  -3.66%  "
Set-TextValue "D34" "404.71"
Set-TextValue "E34" "  -3.49%  "
Set-TextValue "D35" "19.09"
Set-TextValue "E35" "  -0.07%  "
Set-TextValue "E36" "  -2.18%  "
Set-TextValue "E37" "  +0.01%  "
Set-TextValue "D38" "4.33"
Set-TextValue "E38" "  -2.23%  "
Set-TextValue "E39" "  -1.56%  "
Set-TextValue "E40" "  +0.04%  "
Set-TextValue "D41" "39.28"
Set-TextValue "E41" "  -2.99%  "
Set-TextValue "D42" "150.91"
Set-TextValue "E42" "  -1.30%  "
Set-TextValue "E43" "  -1.70%  "
Set-TextValue "E44" "  -0.89%  "
Set-TextValue "D45" "0.0533"
Set-TextValue "E45" "  +0.27%  "
Set-TextValue "E46" "  -1.94%  "
Set-TextValue "E47" "  -0.58%  "
Set-TextValue "E48" "  +2.19%  "
Set-TextValue "D49" "18.20"
Set-TextValue "E49" "  -3.26%  "
Set-TextValue "D50" "11.28"
Set-TextValue "E50" "  +0.23%  "
Set-TextValue "E51" "  -8.06%  "
